$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D1: replace note text with new multi-run rich text ---
$part1 = "*Note :`n- Kolom isian tidak boleh kosong`n- Kolom Gambar Rak, isi dengan default_gambar_rak.jpg`n- Pastikan seluruh kolom yang kosong "
$part2 = "dihapus"
$part3 = ", dengan cara blok kolom yang kosong kemudian klik kanan -> pilih delete  -> ok"
$fullText = $part1 + $part2 + $part3

$d1 = $ws.Cells.Item(1,4)
$d1.Value = $fullText

$start2 = $part1.Length + 1
$len2 = $part2.Length
$d1.Characters($start2, $len2).Font.Bold = $true
$d1.Characters($start2, $len2).Font.Color = 255
$d1.Characters($start2, $len2).Font.Name = "Calibri"
$d1.Characters($start2, $len2).Font.Size = 11

$start3 = $start2 + $len2
$len3 = $part3.Length
$d1.Characters($start3, $len3).Font.Color = 255
$d1.Characters($start3, $len3).Font.Name = "Calibri"
$d1.Characters($start3, $len3).Font.Size = 11

# --- Header row formatting: A1:C1 bold + centered, D1 red + left/wrap ---
$headerRange = $ws.Range("A1:C1")
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108

$d1.HorizontalAlignment = -4131
$d1.VerticalAlignment = -4108
$d1.WrapText = $true

$ws.Rows.Item(1).RowHeight = 105

# --- Hide unused columns beyond D ---
$ws.Range("E1:XFC1").EntireColumn.Hidden = $true

# --- Selection state ---
$ws.Range("A3:D1048576").Select()
